$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.549.41'
$ws.Range('E2').Value = '  +0.97%  '
$ws.Range('D3').Value = '1.889.50'
$ws.Range('E3').Value = '  +1.26%  '
$ws.Range('D4').Value = '0.9996'
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range('D5').Value = '244.57'
$ws.Range('E5').Value = '  +4.21%  '
$ws.Range('D6').Value = '0.9992'
$ws.Range('E6').Value = '  -0.17%  '
$ws.Range('D7').Value = '0.4770'
$ws.Range('E7').Value = '  +1.54%  '
$ws.Range('D8').Value = '0.2901'
$ws.Range('E8').Value = '  +1.83%  '
$ws.Range('D9').Value = '42.82'
$ws.Range('E9').Value = '  +2.37%  '
$ws.Range('D10').Value = '0.06550'
$ws.Range('E10').Value = '  -0.14%  '
$ws.Range('D11').Value = '21.40'
$ws.Range('E11').Value = '  +0.47%  '
$ws.Range('D12').Value = '0.07770'
$ws.Range('E12').Value = '  -0.66%  '
$ws.Range('D13').Value = '1.902.84'
$ws.Range('E13').Value = '  +1.80%  '
$ws.Range('D14').Value = '0.7379'
$ws.Range('E14').Value = '  +6.75%  '
$ws.Range('D15').Value = '96.55'
$ws.Range('E15').Value = '  -0.26%  '
$ws.Range('D16').Value = '5.175'
$ws.Range('E16').Value = '  +1.94%  '
$ws.Range('D17').Value = '276.83'
$ws.Range('E17').Value = '  +3.47%  '
$ws.Range('D18').Value = '30.542.97'
$ws.Range('E18').Value = '  +0.82%  '
$ws.Range('D19').Value = '13.62'
$ws.Range('E19').Value = '  -0.97%  '
$ws.Range('D20').Value = '0.000007613'
$ws.Range('E20').Value = '  -1.26%  '
$ws.Range('D21').Value = '0.9992'
$ws.Range('E21').Value = '  -0.14%  '
$ws.Range('D22').Value = '2.138.77'
$ws.Range('E22').Value = '  +0.86%  '
$ws.Range('D23').Value = '5.312'
$ws.Range('E23').Value = '  +1.03%  '
$ws.Range('D24').Value = '0.9994'
$ws.Range('E24').Value = '  -0.16%  '
$ws.Range('D25').Value = '6.218'
$ws.Range('E25').Value = '  +1.06%  '
$ws.Range('D26').Value = '9.321'
$ws.Range('E26').Value = '  -2.03%  '
$ws.Range('D27').Value = '165.15'
$ws.Range('E27').Value = '  -0.51%  '
$ws.Range('D28').Value = '19.09'
$ws.Range('E28').Value = '  +1.40%  '
$ws.Range('D29').Value = '1.993'
$ws.Range('E29').Value = '  +3.11%  '
$ws.Range('D30').Value = '1.383'
$ws.Range('E30').Value = '  +1.45%  '
$ws.Range('D31').Value = '0.09976'
$ws.Range('E31').Value = '  +0.86%  '
$ws.Range('D32').Value = '1.513'
$ws.Range('E32').Value = '  +4.02%  '
$ws.Range('D33').Value = '4.347'
$ws.Range('E33').Value = '  -0.02%  '
$ws.Range('D34').Value = '4.113'
$ws.Range('E34').Value = '  +1.66%  '
$ws.Range('D35').Value = '0.04781'
$ws.Range('E35').Value = '  +0.89%  '
$ws.Range('E36').Value = '  +0.36%  '
$ws.Range('D37').Value = '0.7040'
$ws.Range('E37').Value = '  +0.17%  '
$ws.Range('E38').Value = '  -0.06%  '
$ws.Range('D39').Value = '0.01856'
$ws.Range('E39').Value = '  -0.58%  '
$ws.Range('D40').Value = '2.767'
$ws.Range('E40').Value = '  -0.08%  '
$ws.Range('D41').Value = '6.511'
$ws.Range('E41').Value = '  +3.20%  '
$ws.Range('D42').Value = '70.89'
$ws.Range('E42').Value = '  -3.31%  '
$ws.Range('D43').Value = '1.932'
$ws.Range('E43').Value = '  -0.74%  '
$ws.Range('B44').Value = 'TheSandbox'
$ws.Range('C44').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D44').Value = '0.4200'
$ws.Range('E44').Value = '  +0.93%  '
$ws.Range('B45').Value = 'TrustWalletToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D45').Value = '0.8454'
$ws.Range('E45').Value = '  +1.50%  '
$ws.Range('D46').Value = '0.9993'
$ws.Range('E46').Value = '  -0.15%  '
$ws.Range('D47').Value = '102.94'
$ws.Range('E47').Value = '  +0.16%  '
$ws.Range('D48').Value = '9.449'
$ws.Range('E48').Value = '  +3.40%  '
$ws.Range('D49').Value = '7.170'
$ws.Range('E49').Value = '  +0.87%  '
$ws.Range('D50').Value = '931.35'
$ws.Range('E50').Value = '  -4.68%  '
$ws.Range('D51').Value = '35.36'
$ws.Range('E51').Value = '  +2.50%  '
